$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6332
$ws.Range("C23").Value = 998
$ws.Range("D23").Value = 5901545
$ws.Range("E23").Value = 932.0191092861655
$ws.Range("F23").Value = 8.647906657515447
$ws.Range("G23").Value = 3.850156087408951
$ws.Range("H23").Value = 26.45500355374097
